# horarios.xlsx edit: rename sheet1 -> Horarios, add Participantes sheet,
# replace email-tuple strings in Horarios with id-tuple strings, and
# populate the Participantes roster.

$wb = $excel.ActiveWorkbook

# --- Rename the existing (only) sheet to "Horarios" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Horarios"

# --- Add a new "Participantes" sheet right after "Horarios" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Participantes"

# --- Update Horarios participant lists (emails -> participant ids) ---
$ws1.Range("B2").Value  = '[''Sin participantes'']'

$ws1.Range("B3").Value  = '[(14,)]'
$ws1.Range("C3").Value  = '[''Sin participantes'']'
$ws1.Range("E3").Value  = '[(14,)]'

$ws1.Range("B4").Value  = '[(14,)]'
$ws1.Range("C4").Value  = '[(14,)]'
$ws1.Range("D4").Value  = '[''Sin participantes'']'
$ws1.Range("E4").Value  = '[(14,)]'

$ws1.Range("B5").Value  = '[(14,)]'
$ws1.Range("C5").Value  = '[(14,)]'
$ws1.Range("D5").Value  = '[(14,)]'
$ws1.Range("E5").Value  = '[''Sin participantes'']'
$ws1.Range("F5").Value  = '[(16,)]'

$ws1.Range("D6").Value  = '[''Sin participantes'']'
$ws1.Range("E6").Value  = '[(14,)]'
$ws1.Range("F6").Value  = '[(14,)]'

$ws1.Range("B7").Value  = '[(15,)]'
$ws1.Range("D7").Value  = '[(15,)]'
$ws1.Range("E7").Value  = '[''Sin participantes'']'
$ws1.Range("F7").Value  = '[(15,)]'

$ws1.Range("B8").Value  = '[(15,), (16,)]'
$ws1.Range("D8").Value  = '[(15,), (14,), (16,)]'
$ws1.Range("F8").Value  = '[(15,), (14,)]'

$ws1.Range("B9").Value  = '[(15,), (16,)]'
$ws1.Range("C9").Value  = '[(14,)]'
$ws1.Range("D9").Value  = '[(15,), (16,)]'
$ws1.Range("F9").Value  = '[(15,)]'

$ws1.Range("B10").Value = '[(16,)]'
$ws1.Range("D10").Value = '[(14,), (16,)]'

$ws1.Range("F11").Value = '[(16,)]'

# --- Populate the Participantes roster sheet ---
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "Nombre"
$ws2.Range("C1").Value = "Apellidos"
$ws2.Range("D1").Value = "Correo"
$ws2.Range("E1").Value = "Teléfono"

# Reuse the bold/bordered/centered header formatting from Horarios!A1:F1
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)

$ws2.Range("A2").Value = 16
$ws2.Range("B2").Value = "JUAN JOSE"
$ws2.Range("C2").Value = "ORTEGA MORALES"
$ws2.Range("D2").Value = "juan.ortega4533@alumnos.udg.mx"
$ws2.Range("E2").Value = "'3322445566"

$ws2.Range("A3").Value = 15
$ws2.Range("B3").Value = "SAMUEL"
$ws2.Range("C3").Value = "CEBALLOS MURGUIA"
$ws2.Range("D3").Value = "samuel.ceballos4453@alumnos.udg.mx"
$ws2.Range("E3").Value = "'3311111111"

$ws2.Range("A4").Value = 14
$ws2.Range("B4").Value = "FRANCO EDUARDO"
$ws2.Range("C4").Value = "SILVA CHACÓN"
$ws2.Range("D4").Value = "franco.silva4477@alumnos.udg.mx"
$ws2.Range("E4").Value = "'3355669988"

# Keep "Horarios" as the active tab (matches original workbook view state)
[void]$ws1.Activate()
[void]$ws1.Range("A1").Select()

Write-Output "ok"
